$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73, shifting existing rows 73:190 down to 74:191.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record's data.
$ws.Range("A73").Value = 9
$ws.Range("B73").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C73").Value = "Metropolitana"
$ws.Range("D73").Value = 44638
$ws.Range("E73").Value = 13
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100101
$ws.Range("H73").Value = "Berries"
$ws.Range("I73").Value = 100101001
$ws.Range("J73").Value = "Arándano (blue)"
$ws.Range("K73").Value = "Sin especificar"
$ws.Range("L73").Value = "Primera"
$ws.Range("M73").Value = 350
$ws.Range("N73").Value = 4000
$ws.Range("O73").Value = 4000
$ws.Range("P73").Value = 4000
$ws.Range("Q73").Value = "`$/bandeja 2 kilos"
$ws.Range("R73").Value = "Provincia de Linares"
$ws.Range("S73").Value = 2000
$ws.Range("T73").Value = 2
